# FAAS-828 Device customization for android
# The discDate test value is updated from "03-04-2021" to "05-12-2021".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdata")

$ws.Range("T2").Value = "05-12-2021"

# Leave the sheet with the edited cell selected/active, matching the
# post-edit view state captured in the workbook.
$ws.Range("T2").Select()
